# Auto-generated: apply cached-value updates from scheduled market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1366.8438
$ws.Range("I17").Value = 900
$ws.Range("J17").Value = 1433.5358
$ws.Range("K17").Value = 2700
$ws.Range("L17").Value = 4300.607400000001
$ws.Range("M17").Value = -2532
$ws.Range("N17").Value = -4636.607400000001
# Row 19
$ws.Range("H19").Value = 2180.0833
$ws.Range("I19").Value = 1432.3334
$ws.Range("J19").Value = 2927.8333
$ws.Range("K19").Value = 1432.3334
$ws.Range("L19").Value = 2927.8333
$ws.Range("M19").Value = -1257.3334
$ws.Range("N19").Value = -3277.8333
# Row 29
$ws.Range("H29").Value = 1122.5
$ws.Range("I29").Value = 625
$ws.Range("J29").Value = 1620
$ws.Range("K29").Value = 1875
$ws.Range("L29").Value = 4860
$ws.Range("M29").Value = -1594
$ws.Range("N29").Value = -5422
# Row 76
$ws.Range("H76").Value = 4122.125
$ws.Range("I76").Value = 4246
$ws.Range("J76").Value = 3998.25
$ws.Range("K76").Value = 4246
$ws.Range("L76").Value = 3998.25
$ws.Range("M76").Value = -3931
$ws.Range("N76").Value = -4628.25
# Row 79
$ws.Range("H79").Value = 4122.125
$ws.Range("I79").Value = 4246
$ws.Range("J79").Value = 3998.25
$ws.Range("K79").Value = 4246
$ws.Range("L79").Value = 3998.25
$ws.Range("M79").Value = -3154
$ws.Range("N79").Value = -6182.25
# Row 137
$ws.Range("H137").Value = 2633.2856
$ws.Range("I137").Value = 2589.95
$ws.Range("K137").Value = 7769.849999999999
$ws.Range("M137").Value = -5219.849999999999

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1548282.6
$ws.Range("I32").Value = 1671556
$ws.Range("J32").Value = 7365.1665
$ws.Range("K32").Value = 1671556
$ws.Range("L32").Value = 7365.1665
$ws.Range("M32").Value = -1671269
$ws.Range("N32").Value = -7939.1665
# Row 132
$ws.Range("H132").Value = 8500.103999999999
$ws.Range("I132").Value = 8777.75
$ws.Range("K132").Value = 26333.25
$ws.Range("M132").Value = -23803.25

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2789.318
$ws.Range("I105").Value = 2084.8
$ws.Range("J105").Value = 4299
$ws.Range("K105").Value = 2084.8
$ws.Range("L105").Value = 4299
$ws.Range("M105").Value = -337.8000000000002
$ws.Range("N105").Value = -7793
# Row 128
$ws.Range("H128").Value = 3391.4443
$ws.Range("I128").Value = 3391.4443
$ws.Range("K128").Value = 10174.3329
$ws.Range("M128").Value = -7684.332900000001

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5892.915
$ws.Range("I31").Value = 2263.2104
$ws.Range("K31").Value = 2263.2104
$ws.Range("M31").Value = -1968.2104
# Row 34
$ws.Range("H34").Value = 5892.915
$ws.Range("I34").Value = 2263.2104
$ws.Range("K34").Value = 2263.2104
$ws.Range("M34").Value = -2061.2104
# Row 58
$ws.Range("H58").Value = 8933206
$ws.Range("I58").Value = 14287255
$ws.Range("J58").Value = 9791.190000000001
$ws.Range("K58").Value = 14287255
$ws.Range("L58").Value = 9791.190000000001
$ws.Range("M58").Value = -14287052
$ws.Range("N58").Value = -10197.19
# Row 62
$ws.Range("H62").Value = 9989
$ws.Range("I62").Value = 9983.75
$ws.Range("K62").Value = 9983.75
$ws.Range("M62").Value = -9359.75
# Row 65
$ws.Range("H65").Value = 9989
$ws.Range("I65").Value = 9983.75
$ws.Range("K65").Value = 49918.75
$ws.Range("M65").Value = -46798.75
# Row 107
$ws.Range("H107").Value = 3076
$ws.Range("I107").Value = 1400
$ws.Range("J107").Value = 3794.2856
$ws.Range("K107").Value = 1400
$ws.Range("L107").Value = 3794.2856
$ws.Range("M107").Value = 520
$ws.Range("N107").Value = -7634.2856
# Row 134
$ws.Range("H134").Value = 9943.842000000001
$ws.Range("I134").Value = 2837.3333
$ws.Range("J134").Value = 11276.3125
$ws.Range("K134").Value = 8511.999899999999
$ws.Range("L134").Value = 33828.9375
$ws.Range("M134").Value = -5976.999899999999
$ws.Range("N134").Value = -38898.9375
# Row 136
$ws.Range("H136").Value = 8933206
$ws.Range("I136").Value = 14287255
$ws.Range("J136").Value = 9791.190000000001
$ws.Range("K136").Value = 42861765
$ws.Range("L136").Value = 29373.57
$ws.Range("M136").Value = -42859215
$ws.Range("N136").Value = -34473.57

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 3333433
$ws.Range("I12").Value = 4
$ws.Range("K12").Value = 12
$ws.Range("M12").Value = 161
# Row 62
$ws.Range("H62").Value = 5604
$ws.Range("I62").Value = 5006
$ws.Range("K62").Value = 15018
$ws.Range("M62").Value = -14332
# Row 65
$ws.Range("H65").Value = 5604
$ws.Range("I65").Value = 5006
$ws.Range("K65").Value = 45054
$ws.Range("M65").Value = -41622
# Row 111
$ws.Range("H111").Value = 22583.334
$ws.Range("I111").Value = 22583.334
$ws.Range("K111").Value = 67750.00199999999
$ws.Range("M111").Value = -64683.00199999999
# Row 138
$ws.Range("H138").Value = 5158.5884
$ws.Range("I138").Value = 4069.077
$ws.Range("J138").Value = 8699.5
$ws.Range("K138").Value = 12207.231
$ws.Range("L138").Value = 26098.5
$ws.Range("M138").Value = -7067.231
$ws.Range("N138").Value = -36378.5

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2259.2083
$ws.Range("I80").Value = 1868.9166
$ws.Range("J80").Value = 2649.5
$ws.Range("K80").Value = 1868.9166
$ws.Range("L80").Value = 2649.5
$ws.Range("M80").Value = -870.9166
$ws.Range("N80").Value = -4645.5
# Row 83
$ws.Range("H83").Value = 2259.2083
$ws.Range("I83").Value = 1868.9166
$ws.Range("J83").Value = 2649.5
$ws.Range("K83").Value = 9344.583000000001
$ws.Range("L83").Value = 13247.5
$ws.Range("M83").Value = -4352.583000000001
$ws.Range("N83").Value = -23231.5
# Row 102
$ws.Range("H102").Value = 2838.9614
$ws.Range("I102").Value = 2627.2273
$ws.Range("K102").Value = 2627.2273
$ws.Range("M102").Value = -1005.2273
# Row 113
$ws.Range("H113").Value = 6807.8857
$ws.Range("I113").Value = 3232.4167
$ws.Range("J113").Value = 8673.348
$ws.Range("K113").Value = 3232.4167
$ws.Range("L113").Value = 8673.348
$ws.Range("M113").Value = -1062.4167
$ws.Range("N113").Value = -13013.348
# Row 132
$ws.Range("H132").Value = 8696.706
$ws.Range("I132").Value = 3209.25
$ws.Range("K132").Value = 9627.75
$ws.Range("M132").Value = -7097.75

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 986.7778
$ws.Range("I16").Value = 1068.4286
$ws.Range("J16").Value = 701
$ws.Range("K16").Value = 1068.4286
$ws.Range("L16").Value = 701
$ws.Range("M16").Value = -898.4286
$ws.Range("N16").Value = -1041
# Row 68
$ws.Range("H68").Value = 1995.8334
$ws.Range("J68").Value = 1995
$ws.Range("L68").Value = 1995
$ws.Range("N68").Value = -3493
# Row 71
$ws.Range("H71").Value = 1995.8334
$ws.Range("J71").Value = 1995
$ws.Range("L71").Value = 9975
$ws.Range("N71").Value = -17463

$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Range("H45").Value = 11222.167
$ws.Range("J45").Value = 12677.2
$ws.Range("L45").Value = 12677.2
$ws.Range("N45").Value = -13659.2
# Row 100
$ws.Range("H100").Value = 1091.9166
$ws.Range("I100").Value = 972.8
$ws.Range("J100").Value = 1177
$ws.Range("K100").Value = 1945.6
$ws.Range("L100").Value = 2354
$ws.Range("M100").Value = -1404.6
$ws.Range("N100").Value = -3436
# Row 122
$ws.Range("H122").Value = 4972.7075
$ws.Range("I122").Value = 4687.3447
$ws.Range("J122").Value = 5662.3335
$ws.Range("K122").Value = 14062.0341
$ws.Range("L122").Value = 16987.0005
$ws.Range("M122").Value = -11612.0341
$ws.Range("N122").Value = -21887.0005
# Row 132
$ws.Range("H132").Value = 9623810
$ws.Range("I132").Value = 11367029
$ws.Range("K132").Value = 34101087
$ws.Range("M132").Value = -34098557
